# Updates the cryptos list (prices + volume%) as produced by the
# "Updated cryptos list ... with GitHub Actions" workflow run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as text so values such as "1.0000",
    # "29.135.68" or "0.000008169" are not reinterpreted by Excel as
    # numbers/dates. ClearFormats() afterwards drops the temporary "@"
    # number-format style again so the cell keeps its original (default)
    # style index.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Rows whose Price (D) and Volume(1h) (E) changed -----------------------

$deChanges = @{
    2  = @("29.135.68", "  -0.23%  ")
    3  = @("1.842.37", "  -0.42%  ")
    4  = @("0.9993", "  -0.03%  ")
    5  = @("241.69", "  -1.74%  ")
    6  = @($null, "  -1.67%  ")
    7  = @($null, "  -0.04%  ")
    8  = @("0.3022", "  -1.31%  ")
    9  = @("0.07476", "  -3.13%  ")
    10 = @("23.16", "  -1.44%  ")
    11 = @("0.07668", "  -1.97%  ")
    12 = @("1.846.55", "  -0.01%  ")
    13 = @($null, "  -1.15%  ")
    14 = @("0.6844", "  -0.18%  ")
    15 = @("87.65", "  -5.76%  ")
    16 = @("6.178", "  -6.85%  ")
    17 = @("29.136.43", "  -0.17%  ")
    18 = @("0.000008169", "  -1.67%  ")
    19 = @("2.080.43", "  -0.34%  ")
    20 = @("228.37", "  -5.46%  ")
    21 = @($null, "  -1.41%  ")
    23 = @("7.416", "  -1.38%  ")
    24 = @("1.0000", "  -0.01%  ")
    25 = @("0.1456", "  -3.73%  ")
    26 = @("160.07", "  +0.42%  ")
    27 = @("8.763", "  -0.73%  ")
    28 = @("18.09", "  -1.05%  ")
    29 = @("1.512", "  -1.81%  ")
    30 = @("4.275", "  +1.15%  ")
    31 = @("4.145", "  -0.81%  ")
    32 = @($null, "  -0.15%  ")
    33 = @("0.05228", "  +2.09%  ")
    34 = @("0.7661", "  -3.45%  ")
    35 = @("1.852", "  -0.77%  ")
    36 = @("1.137", "  -0.93%  ")
    37 = @("2.680", "  -0.53%  ")
    38 = @("1.314.14", "  -0.44%  ")
    39 = @("0.01839", "  -1.85%  ")
    40 = @("2.728", "  +0.59%  ")
    41 = @("0.9351", "  -1.42%  ")
    42 = @("105.13", "  -1.98%  ")
    43 = @("5.798", "  -3.35%  ")
    44 = @("0.9993", "  -0.12%  ")
    50 = @("1.776", "  +0.58%  ")
    51 = @("0.07513", "  +19.12%  ")
}

foreach ($r in $deChanges.Keys) {
    $vals = $deChanges[$r]
    $price = $vals[0]
    $volume = $vals[1]
    if ($null -ne $price) {
        Set-TextValue $ws.Range("D$r") $price
    }
    Set-TextValue $ws.Range("E$r") $volume
}

# --- Rows 45-49: the coin list was re-ordered, shifting BabyDogeCoin to the
# top of this block and pushing RocketPoolETH / Mantle / Aave / EnergySwap
# down by one row. Coin name (B), Link (C), Price (D) and Volume(1h) (E)
# all change for these rows.

$rows45to49 = @{
    45 = @("BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.00000000123", "  +3.80%  ")
    46 = @("RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "1.983.76", "  -0.25%  ")
    47 = @("Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.5198", "  +0.33%  ")
    48 = @("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "64.88", "  +1.19%  ")
    49 = @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.501", "  -2.21%  ")
}

foreach ($r in $rows45to49.Keys) {
    $vals = $rows45to49[$r]
    Set-TextValue $ws.Range("B$r") $vals[0]
    Set-TextValue $ws.Range("C$r") $vals[1]
    Set-TextValue $ws.Range("D$r") $vals[2]
    Set-TextValue $ws.Range("E$r") $vals[3]
}
